# Update "F想去人数" (want-to-go count) values on the "展览" and "全部类型"
# sheets to reflect newly scraped numbers (gh-pages output regenerated at 456a3b4).

$wb = $excel.ActiveWorkbook

# Row -> new F-column value for sheet "展览"
$exhibitionUpdates = @{
    2  = 1345
    3  = 1212
    4  = 902
    11 = 2425
    12 = 1601
    13 = 1455
    15 = 245
    16 = 596
    17 = 778
    18 = 65
    19 = 306
    22 = 26
    23 = 530
    24 = 4929
    26 = 456
    27 = 75
    28 = 158
    29 = 136
    30 = 222
    31 = 97
    32 = 27
    34 = 711
    36 = 50
    38 = 389
    39 = 1040
    42 = 167
}

# Row -> new F-column value for sheet "全部类型"
$allTypesUpdates = @{
    2  = 1345
    5  = 1212
    6  = 902
    17 = 2425
    18 = 1601
    19 = 1455
    21 = 245
    22 = 596
    24 = 778
    25 = 65
    26 = 306
    28 = 26
    29 = 530
    30 = 4929
    32 = 456
    33 = 75
    34 = 158
    35 = 136
    36 = 222
    37 = 97
    38 = 27
    40 = 711
    41 = 50
    42 = 389
    43 = 1040
    45 = 167
}

$wsExhibition = $wb.Worksheets.Item("展览")
foreach ($row in $exhibitionUpdates.Keys) {
    $wsExhibition.Range("F$row").Value = $exhibitionUpdates[$row]
}

$wsAllTypes = $wb.Worksheets.Item("全部类型")
foreach ($row in $allTypesUpdates.Keys) {
    $wsAllTypes.Range("F$row").Value = $allTypesUpdates[$row]
}
